$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.954.95'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '1.651.07'
$ws.Range("E3").Value = '  +2.85%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.86'
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  +2.30%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("E9").Value = '  +1.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.19'
$ws.Range("E10").Value = '  +5.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0880'
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").Value = '1.884.20'
$ws.Range("E12").Value = '  +2.89%  '
$ws.Range("D13").Value = '1.658.75'
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.520'
$ws.Range("E15").Value = '  +2.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.11'
$ws.Range("E16").Value = '  +2.89%  '
$ws.Range("D17").Value = '26.952.15'
$ws.Range("E17").Value = '  +2.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '235.64'
$ws.Range("E18").Value = '  +2.05%  '
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.75'
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  +3.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("E23").Value = '  +4.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  +2.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.37'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.81'
$ws.Range("E29").Value = '  +2.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("D32").Value = '1.562.95'
$ws.Range("E32").Value = '  +4.67%  '
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("E34").Value = '  +5.06%  '
$ws.Range("E35").Value = '  +9.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("E37").Value = '  +4.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.896'
$ws.Range("E38").Value = '  +9.39%  '
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("E40").Value = '  +3.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.04'
$ws.Range("E42").Value = '  +8.58%  '
$ws.Range("E43").Value = '  +2.44%  '
$ws.Range("D44").Value = '1.791.72'
$ws.Range("E44").Value = '  +2.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.777'
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.92'
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("E48").Value = '  +1.96%  '
$ws.Range("E49").Value = '  +3.17%  '
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("E51").Value = '  +2.95%  '
